$d = $word.ActiveDocument

# Append three new bulleted list items after the last paragraph, matching
# the existing "ListParagraph" / numId=2 bullet list used above it.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Ukicer 24"

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Last
$p1.Range.InsertParagraphAfter()

$d = $word.ActiveDocument
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Mahara " + [char]0x2013 + " reflective learning in moodle"

$d = $word.ActiveDocument
$p2 = $d.Paragraphs.Last
$p2.Range.InsertParagraphAfter()

$d = $word.ActiveDocument
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Jupyter physics server"
